# Apply the "conclusion" tutorial deck edit:
#  1. Delete 3 slides (old slide 2 "Arkouda Performance Results", old slide 3
#     "Arachne Performance Results" with images, old slide 4 "Towards Arachne
#     1.5 and Beyond").
#  2. Tweak bold formatting in the author line on the title slide.
#  3. Retitle / recolor a few remaining slides ("Arachne" -> red, "Arkouda" ->
#     purple) and rewrite the conclusion bullets.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Remove the three slides that were deleted from the deck.
#    (Positions 2, 3, 4 in the original ordering -- deleting position 2
#    three times removes exactly those three slides because everything
#    shifts down.)
# ---------------------------------------------------------------------------
$p.Slides.Item(2).Delete()
$p.Slides.Item(2).Delete()
$p.Slides.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2. Title slide (slide 1): rebalance bold runs in the author list.
#    Before: **Oliver Alvarado Rodriguez**, Naren Khatwani, **Zhihui** Du, David Bader
#    After : Oliver Alvarado Rodriguez, Naren Khatwani, **Zhihui Du**, David Bader
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2).TextFrame.TextRange
$name1 = $subtitle.Characters(13, 27)   # "Oliver Alvarado Rodriguez, "
$name1.Text = "Oliver Alvarado Rodriguez, "
$name1.Font.Bold = $false
$subtitle.Characters(56, 9).Font.Bold = $true      # "Zhihui Du"

# ---------------------------------------------------------------------------
# 3. New slide 2 (was "Hybrid Arachne"): retitle + recolor "Arachne", resize
#    title box, recolor the "Arachne 1.5" pill label.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1)
$title2.Left = 36
$title2.Top = 18
$title2.Width = 858
$title2.Height = 104.37503937007874
$title2.TextFrame.TextRange.Text = "Scaling Arachne from SMPs to MPPs & Clusters"
$title2.TextFrame.TextRange.Characters(9, 7).Font.Color.RGB = 255

$pill2 = $s2.Shapes.Item(28)
$pill2.TextFrame.TextRange.Text = "Arachne 1.5"
$pill2.TextFrame.TextRange.Characters(1, 7).Font.Color.RGB = 255

# ---------------------------------------------------------------------------
# 4. New slide 3 (was "Enhancing Arachne for Property Graphs"): retitle with
#    "(2.0)" and recolor "Arachne".
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "Enhancing Arachne (2.0) for Property Graphs"
$title3.TextFrame.TextRange.Characters(11, 7).Font.Color.RGB = 255

# ---------------------------------------------------------------------------
# 5. New slide 4 ("Conclusion"): rewrite the bullet list and recolor the
#    product names.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "We have shown the usability of Arkouda for large-scale data analysis.`rWe have shown proof of concept of Arachne through breadth-first search, truss analytics, connected components, etc. `rWe have outlined our goals of fleshing out Arachne to be a hybrid solution for (property) graph analysis scaling from SMPs to MPPs and clusters`rWe have outlined the blueprint for the future of Arachne."
$body4.Characters(32, 7).Font.Color.RGB = 10498160   # "Arkouda" -> purple 7030A0
$body4.Characters(105, 7).Font.Color.RGB = 255        # "Arachne" (2nd paragraph) -> red FF0000
$body4.Characters(231, 7).Font.Color.RGB = 255        # "Arachne" (3rd paragraph) -> red FF0000
$body4.Characters(381, 7).Font.Color.RGB = 255        # "Arachne" (4th paragraph) -> red FF0000
